$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.309.93'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '1.928.09'

$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '0.7492'
$ws.Range('E5').Value = '  +4.96%  '

$ws.Range('D6').Value = '243.51'
$ws.Range('E6').Value = '  -3.10%  '

$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('D8').Value = '0.3156'

$ws.Range('D9').Value = '27.46'
$ws.Range('E9').Value = '  +0.26%  '

$ws.Range('D10').Value = '0.06973'
$ws.Range('E10').Value = '  -3.11%  '

$ws.Range('D11').Value = '0.08005'
$ws.Range('E11').Value = '  -1.09%  '

$ws.Range('D12').Value = '0.7700'
$ws.Range('E12').Value = '  -3.69%  '

$ws.Range('D13').Value = '1.931.09'
$ws.Range('E13').Value = '  -0.03%  '

$ws.Range('D14').Value = '5.323'
$ws.Range('E14').Value = '  -2.05%  '

$ws.Range('D15').Value = '93.20'
$ws.Range('E15').Value = '  -1.68%  '

$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '14.33'
$ws.Range('E16').Value = '  -3.41%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '30.296.38'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('D18').Value = '250.68'
$ws.Range('E18').Value = '  -1.02%  '

$ws.Range('D19').Value = '0.000007892'
$ws.Range('E19').Value = '  -2.59%  '

$ws.Range('D20').Value = '5.748'
$ws.Range('E20').Value = '  -0.98%  '

$ws.Range('D21').Value = '2.186.91'
$ws.Range('E21').Value = '  +0.25%  '

$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = '6.631'
$ws.Range('E24').Value = '  -4.16%  '

$ws.Range('D25').Value = '9.421'
$ws.Range('E25').Value = '  -2.96%  '

$ws.Range('D26').Value = '165.84'
$ws.Range('E26').Value = '  +0.63%  '

$ws.Range('D27').Value = '18.90'
$ws.Range('E27').Value = '  -1.87%  '

$ws.Range('D28').Value = '0.1321'
$ws.Range('E28').Value = '  +3.18%  '

$ws.Range('E29').Value = '  -5.62%  '

$ws.Range('E30').Value = '  +0.77%  '

$ws.Range('D31').Value = '1.510'
$ws.Range('E31').Value = '  -2.12%  '

$ws.Range('D32').Value = '4.368'
$ws.Range('E32').Value = '  -1.46%  '

$ws.Range('D33').Value = '4.098'
$ws.Range('E33').Value = '  -2.62%  '

$ws.Range('D34').Value = '0.05098'
$ws.Range('E34').Value = '  -2.13%  '

$ws.Range('D35').Value = '1.279'
$ws.Range('E35').Value = '  +0.89%  '

$ws.Range('D36').Value = '0.7432'
$ws.Range('E36').Value = '  -0.95%  '

$ws.Range('D37').Value = '2.778'
$ws.Range('E37').Value = '  +0.35%  '

$ws.Range('D38').Value = '0.01946'
$ws.Range('E38').Value = '  -1.01%  '

$ws.Range('D39').Value = '2.796'
$ws.Range('E39').Value = '  -0.20%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '6.414'
$ws.Range('E40').Value = '  -0.36%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '76.97'
$ws.Range('E41').Value = '  -2.63%  '

$ws.Range('D42').Value = '0.4429'
$ws.Range('E42').Value = '  -2.27%  '

$ws.Range('E43').Value = '  -3.77%  '

$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.07%  '

$ws.Range('D45').Value = '0.8314'
$ws.Range('E45').Value = '  -1.30%  '

$ws.Range('D46').Value = '100.76'
$ws.Range('E46').Value = '  -1.14%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '7.435'
$ws.Range('E47').Value = '  -0.23%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.664'
$ws.Range('E48').Value = '  -1.70%  '

$ws.Range('D49').Value = '37.16'

$ws.Range('D50').Value = '974.69'
$ws.Range('E50').Value = '  +9.62%  '

$ws.Range('D51').Value = '0.06032'
$ws.Range('E51').Value = '  -1.04%  '
